$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row before row 35, shifting existing rows 35..114 down to 36..115.
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row's September details/date columns (R/S)
# with the new transaction entry.
$ws.Cells.Item(35, 18).Value = "transfer freedom share anyone axis"
$ws.Cells.Item(35, 19).Value = "2024-09-09 11:56:19"
